$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.048.85"
$ws.Range("E2").Value = "  +0.21%  "

$ws.Range("D3").Value = "1.820.96"
$ws.Range("E3").Value = "  +0.05%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.03"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.27%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.620"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.04%  "

$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "39.83"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -4.86%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.324"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +4.77%  "

$ws.Range("E10").Value = "  -0.18%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0993"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.92%  "

$ws.Range("D12").Value = "2.084.21"
$ws.Range("E12").Value = "  +0.02%  "

$ws.Range("E13").Value = "  +2.34%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.841.93"
$ws.Range("E14").Value = "  +0.79%  "

$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.668"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.33%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.65"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.16%  "

$ws.Range("D17").Value = "35.048.13"
$ws.Range("E17").Value = "  +0.38%  "

$ws.Range("E18").Value = "  +0.30%  "

$ws.Range("E19").Value = "  +0.17%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "240.86"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.08%  "

$ws.Range("E21").Value = "  +2.65%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.70"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.70%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.28"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.75%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "174.00"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.16%  "

$ws.Range("E26").Value = "  +0.48%  "

$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.124"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.00%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.38"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.19%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.53"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.27%  "

$ws.Range("E30").Value = "  +0.12%  "

$ws.Range("E31").Value = "  +2.94%  "

$ws.Range("E32").Value = "  +0.26%  "

$ws.Range("E33").Value = "  -0.39%  "

$ws.Range("E34").Value = "  +12.10%  "

$ws.Range("E35").Value = "  +3.45%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.699"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.47%  "

$ws.Range("E37").Value = "  +0.17%  "

$ws.Range("E38").Value = "  +7.36%  "

$ws.Range("D39").Value = "1.341.48"
$ws.Range("E39").Value = "  +2.12%  "

$ws.Range("E40").Value = "  +1.11%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.989"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.36%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.75"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.56%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.28"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.01%  "

$ws.Range("E44").Value = "  -0.83%  "

$ws.Range("E45").Value = "  -0.20%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.26"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0520"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.05%  "

$ws.Range("E48").Value = "  +0.11%  "

$ws.Range("E49").Value = "  +0.12%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0667"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +4.40%  "

$ws.Range("E51").Value = "  +12.44%  "

